# Atualizacao do documento de homologacao
# Reproduces the commit: the four "Ok" / "nok" status rows (18,19,21,22) of the
# "MS e INC" item are resolved: STATUS becomes "Corrigido" (green-filled) and
# the old RESPOSTA "Ok" note is removed from column F; the RESPOSTA text for
# the other item (rows 20/24/25) is updated to explain what was told to Vini.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Drop the four "Ok" notes in column F for rows 18, 19, 21, 22 -- these
#    cells are removed entirely (not just blanked), matching how Excel drops
#    a cell node once it is fully cleared.
$ws.Range("F18").Clear()
$ws.Range("F19").Clear()
$ws.Range("F21").Clear()
$ws.Range("F22").Clear()

# 2) Mark the corresponding STATUS cells (column C) as resolved, and paint
#    them green to flag the fix -- same green used by the existing
#    conditional-formatting "Ok" rule (FF92D050), applied directly this time.
foreach ($addr in @("C18", "C19", "C21", "C22")) {
    $cell = $ws.Range($addr)
    $cell.Value = "Corrigido"
    $cell.Interior.Color = 5296274
}

# 3) Update the RESPOSTA text for the other pending item (row 20) and extend
#    the same note down into rows 24 and 25, which previously had no note.
$ws.Range("F20").Value = "Expliquei pro Vini o que precisa ser feito."
$ws.Range("F24").Value = "Expliquei pro Vini o que precisa ser feito."
$ws.Range("F25").Value = "Expliquei pro Vini o que precisa ser feito."

# 4) Move the on-screen selection down to the area just edited.
$ws.Range("F25").Select()
